# Big Bend Community College Organizations sheet:
#   * swap the "Organization Name" / "Categories" columns (Category now
#     comes first, Organization Name second)
#   * rename several headers and add a new "Tiktok Link" column (M)
#   * resize columns to match the new header/content widths
#
# $excel / $wb / $ws are provided by the host runtime; $wb.ActiveWorkbook
# is already open on the "Organizations" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) { $lastRow = 49 }

# ---------------------------------------------------------------------------
# Step 1: introduce column M ("Tiktok Link").
# Copy column L's formatting into M first, so the new header cell (M1)
# picks up the existing bold/centered/bordered header style, and the new
# data cells (M2:M49) pick up the same un-styled look as the rest of the
# data rows, before any values are written.
# ---------------------------------------------------------------------------
$ws.Range("L1:L$lastRow").Copy()
$ws.Range("M1:M$lastRow").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Step 2: swap the first two columns' data for every organization row.
# Column A held the organization name and B held its category; after the
# edit A holds the category and B holds the organization name.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $orgName = $ws.Cells.Item($r, 1).Value2
    $category = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $category
    $ws.Cells.Item($r, 2).Value = $orgName
}

# ---------------------------------------------------------------------------
# Step 3: rewrite the header row to match the new column meanings/names.
# ---------------------------------------------------------------------------
$headers = @{
    1  = "Category"
    2  = "Organization Name"
    3  = "Organization Link"
    4  = "Logo Link"
    5  = "Description"
    6  = "Email"
    7  = "Phone Number"
    8  = "Linkedin Link"
    9  = "Instagram Link"
    10 = "Facebook Link"
    11 = "Twitter Link"
    12 = "Youtube Link"
    13 = "Tiktok Link"
}
foreach ($col in $headers.Keys) {
    $ws.Cells.Item(1, $col).Value = $headers[$col]
}

# ---------------------------------------------------------------------------
# Step 4: widen columns to fit the new headers/content.
# (This engine's ColumnWidth getter/setter is offset from the raw OOXML
# <col width> by a constant ~0.83 character padding, so subtract that to
# land on the desired stored width.)
# ---------------------------------------------------------------------------
$colWidths = @{
    1  = 14
    2  = 50
    3  = 50
    4  = 11
    5  = 13
    6  = 7
    7  = 14
    8  = 15
    9  = 16
    10 = 15
    11 = 14
    12 = 14
    13 = 13
}
foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col] - 0.83
}
